$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append the new FRED observations that came in ---
$data = $wb.Worksheets.Item("Data")

# Copy the formatting of the last existing date cell onto the new rows
# before writing their values, so the appended rows look exactly like the
# rest of the series (centered, bordered, date-formatted column A).
$data.Range("A459").Copy($data.Range("A460"))
$data.Range("A459").Copy($data.Range("A461"))
$data.Range("A459").Copy($data.Range("A462"))

$data.Cells.Item(460, 1).Value = 45142
$data.Cells.Item(460, 2).Value = 1793.804

$data.Cells.Item(461, 1).Value = 45145
$data.Cells.Item(461, 2).Value = 1810.583

$data.Cells.Item(462, 1).Value = 45146
$data.Cells.Item(462, 2).Value = 1778.351

# --- Sheet "SeriesInfo": refresh the FRED metadata timestamps ---
$info = $wb.Worksheets.Item("SeriesInfo")

# Format as text first so these date-shaped strings are stored literally
# instead of being auto-converted to date serials by Excel.
$textCells = $info.Range("B3:B4")
$textCells.NumberFormat = "@"
$info.Range("B3").Value = "2023-08-09"
$info.Range("B4").Value = "2023-08-09"

$info.Range("B7").NumberFormat = "@"
$info.Range("B7").Value = "2023-08-08"

# This one already includes a UTC offset suffix, which Excel can't parse as
# a date/time on its own, so it's kept as plain text automatically.
$info.Range("B14").Value = "2023-08-08 13:01:06-05"
